$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3745.4443
$ws.Range("I76").Value = 3281.2
$ws.Range("J76").Value = 4325.75
$ws.Range("K76").Value = 3281.2
$ws.Range("L76").Value = 4325.75
$ws.Range("M76").Value = -2966.2
$ws.Range("N76").Value = -4955.75
$ws.Range("H79").Value = 3745.4443
$ws.Range("I79").Value = 3281.2
$ws.Range("J79").Value = 4325.75
$ws.Range("K79").Value = 3281.2
$ws.Range("L79").Value = 4325.75
$ws.Range("M79").Value = -2189.2
$ws.Range("N79").Value = -6509.75
$ws.Range("H92").Value = 760.85
$ws.Range("I92").Value = 662.9231
$ws.Range("J92").Value = 942.7143
$ws.Range("K92").Value = 662.9231
$ws.Range("L92").Value = 942.7143
$ws.Range("M92").Value = 585.0769
$ws.Range("N92").Value = -3438.7143
$ws.Range("H125").Value = 1084
$ws.Range("I125").Value = 368
$ws.Range("J125").Value = 1322.6666
$ws.Range("K125").Value = 3312
$ws.Range("L125").Value = 11903.9994
$ws.Range("M125").Value = -852
$ws.Range("N125").Value = -16823.9994
$ws.Range("H138").Value = 3643.1829
$ws.Range("I138").Value = 2058.0334
$ws.Range("J138").Value = 4398.016
$ws.Range("K138").Value = 6174.100199999999
$ws.Range("L138").Value = 13194.048
$ws.Range("M138").Value = -1034.100199999999
$ws.Range("N138").Value = -23474.048
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10816.75
$ws.Range("I32").Value = 8955.77
$ws.Range("J32").Value = 29633.334
$ws.Range("K32").Value = 8955.77
$ws.Range("L32").Value = 29633.334
$ws.Range("M32").Value = -8668.77
$ws.Range("N32").Value = -30207.334
$ws.Range("H61").Value = 1696.8684
$ws.Range("I61").Value = 1374.3214
$ws.Range("J61").Value = 2600
$ws.Range("K61").Value = 1374.3214
$ws.Range("L61").Value = 2600
$ws.Range("M61").Value = -1162.3214
$ws.Range("N61").Value = -3024
$ws.Range("H97").Value = 3005.2354
$ws.Range("I97").Value = 4205.8
$ws.Range("J97").Value = 1290.1428
$ws.Range("K97").Value = 4205.8
$ws.Range("L97").Value = 1290.1428
$ws.Range("M97").Value = -3709.8
$ws.Range("N97").Value = -2282.1428
$ws.Range("H136").Value = 1696.8684
$ws.Range("I136").Value = 1374.3214
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 4122.9642
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -1572.9642
$ws.Range("N136").Value = -12900
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4859.102
$ws.Range("I94").Value = 765.8095
$ws.Range("K94").Value = 765.8095
$ws.Range("M94").Value = -314.8095
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2678.3076
$ws.Range("I31").Value = 1265.5652
$ws.Range("J31").Value = 3798.7585
$ws.Range("K31").Value = 1265.5652
$ws.Range("L31").Value = 3798.7585
$ws.Range("M31").Value = -970.5652
$ws.Range("N31").Value = -4388.7585
$ws.Range("H34").Value = 2678.3076
$ws.Range("I34").Value = 1265.5652
$ws.Range("J34").Value = 3798.7585
$ws.Range("K34").Value = 1265.5652
$ws.Range("L34").Value = 3798.7585
$ws.Range("M34").Value = -1063.5652
$ws.Range("N34").Value = -4202.7585
$ws.Range("H132").Value = 1488.0333
$ws.Range("I132").Value = 1169.72
$ws.Range("J132").Value = 3079.6
$ws.Range("K132").Value = 3509.16
$ws.Range("L132").Value = 9238.799999999999
$ws.Range("M132").Value = -979.1599999999999
$ws.Range("N132").Value = -14298.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 3085.5908
$ws.Range("I51").Value = 600
$ws.Range("K51").Value = 1800
$ws.Range("M51").Value = -1340
$ws.Range("H112").Value = 5775.773
$ws.Range("I112").Value = 2850.8
$ws.Range("J112").Value = 6150.769
$ws.Range("K112").Value = 8552.400000000001
$ws.Range("L112").Value = 18452.307
$ws.Range("M112").Value = -7444.400000000001
$ws.Range("N112").Value = -20668.307
$ws.Range("H140").Value = 1389.697
$ws.Range("I140").Value = 1398.0952
$ws.Range("K140").Value = 4194.2856
$ws.Range("M140").Value = 985.7143999999998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2408.1538
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 2350.75
$ws.Range("K80").Value = 2500
$ws.Range("L80").Value = 2350.75
$ws.Range("M80").Value = -1502
$ws.Range("N80").Value = -4346.75
$ws.Range("H83").Value = 2408.1538
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 2350.75
$ws.Range("K83").Value = 12500
$ws.Range("L83").Value = 11753.75
$ws.Range("M83").Value = -7508
$ws.Range("N83").Value = -21737.75
$ws.Range("H97").Value = 3186.6667
$ws.Range("I97").Value = 1763.6364
$ws.Range("J97").Value = 7100
$ws.Range("K97").Value = 1763.6364
$ws.Range("L97").Value = 7100
$ws.Range("M97").Value = -1267.6364
$ws.Range("N97").Value = -8092
$ws.Range("H102").Value = 1393.3846
$ws.Range("I102").Value = 1471
$ws.Range("J102").Value = 1134.6666
$ws.Range("K102").Value = 1471
$ws.Range("L102").Value = 1134.6666
$ws.Range("M102").Value = 151
$ws.Range("N102").Value = -4378.6666
$ws.Range("H113").Value = 2593.2
$ws.Range("I113").Value = 1226.6
$ws.Range("J113").Value = 3959.8
$ws.Range("K113").Value = 1226.6
$ws.Range("L113").Value = 3959.8
$ws.Range("M113").Value = 943.4000000000001
$ws.Range("N113").Value = -8299.799999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1849.2
$ws.Range("I93").Value = 1976.2222
$ws.Range("J93").Value = 1745.2727
$ws.Range("K93").Value = 1976.2222
$ws.Range("L93").Value = 1745.2727
$ws.Range("M93").Value = -728.2221999999999
$ws.Range("N93").Value = -4241.2727
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 83336240
$ws.Range("I81").Value = 250003420
$ws.Range("J81").Value = 2650
$ws.Range("K81").Value = 500006840
$ws.Range("L81").Value = 5300
$ws.Range("M81").Value = -500005779
$ws.Range("N81").Value = -7422
$ws.Range("H84").Value = 83336240
$ws.Range("I84").Value = 250003420
$ws.Range("J84").Value = 2650
$ws.Range("K84").Value = 2500034200
$ws.Range("L84").Value = 26500
$ws.Range("M84").Value = -2500028896
$ws.Range("N84").Value = -37108
$ws.Range("H100").Value = 5311.136
$ws.Range("I100").Value = 12606.875
$ws.Range("J100").Value = 1142.1428
$ws.Range("K100").Value = 25213.75
$ws.Range("L100").Value = 2284.2856
$ws.Range("M100").Value = -24672.75
$ws.Range("N100").Value = -3366.2856
$ws.Range("H136").Value = 5225.1787
$ws.Range("I136").Value = 1591.1765
$ws.Range("J136").Value = 10841.363
$ws.Range("K136").Value = 4773.529500000001
$ws.Range("L136").Value = 32524.089
$ws.Range("M136").Value = -2223.529500000001
$ws.Range("N136").Value = -37624.089
